$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-StatCell($addr, $value, $colorHex) {
    $c = $ws.Range($addr)
    $c.Font.Color = 0x333333
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4108
    $c.Interior.Color = $colorHex
    $c.Value = $value
}

# --- Update existing rows 2-4 (cloud coverage stats changed) ---

# Row 2: O2 unchanged style (green), Q2 4 -> 0 value + fill changes to blue
$ws.Range("O2").Value = 0
Set-StatCell "Q2" 0 0xC88C55

# Row 3: O3 -> orange, P3 -> light blue, Q3 -> medium blue
Set-StatCell "O3" 25 0x70B0F0
Set-StatCell "P3" 25 0xD6A980
Set-StatCell "Q3" 11 0xCD9766

# Row 4: O4 -> red (new style slot, same color), P4 -> very light (new style slot, same color)
Set-StatCell "O4" 94 0x7F70F0
Set-StatCell "P4" 92 0xF9F3EE

# --- Add new row 5 for 29.12.2025 ---
$ws.Range("A5").Value = "29.12.2025"
$ws.Range("B5").Value = 20
$ws.Range("C5").Value = "05:21"
$ws.Range("D5").Value = "10:27"
$ws.Range("E5").Value = "05:05:49"
$ws.Range("F5").Value = "05:08:22"
$ws.Range("G5").Value = "05:11:02"
$ws.Range("H5").Value = "05:13:43"
$ws.Range("I5").Value = "05:16:16"
$ws.Range("J5").Value = "7°"
$ws.Range("K5").Value = "05:07:45"
$ws.Range("L5").Value = -17.2
$ws.Range("M5").Value = "B"
$ws.Range("N5").Value = "3"

Set-StatCell "O5" 33 0x70B0F0
Set-StatCell "P5" 33 0xDBB490
Set-StatCell "Q5" 1 0xC88C55
Set-StatCell "R5" 2 0xC88C55

# --- Extend conditional formatting ranges to include row 5 ---
$lcf = $ws.Range("L2:L4").FormatConditions
$lcf.Item(1).ModifyAppliesToRange($ws.Range("L2:L5"))

$ncf = $ws.Range("N2:N4").FormatConditions
$ncf.Item(1).ModifyAppliesToRange($ws.Range("N2:N5"))
